$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.883.17"
$ws.Range("E2").Value = "  +2.60%  "

$ws.Range("D3").Value = "1.871.99"
$ws.Range("E3").Value = "  +0.94%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  -0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.80"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4842"
$ws.Range("E7").Value = "  +1.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3829"
$ws.Range("E8").Value = "  +3.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07373"
$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9405"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.04"
$ws.Range("E11").Value = "  +5.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07818"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "1.881.49"
$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.498"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.618"
$ws.Range("E15").Value = "  +1.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.25"
$ws.Range("E16").Value = "  +1.74%  "

$ws.Range("E17").Value = "  -0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008882"
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("E19").Value = "  -0.55%  "

$ws.Range("D20").Value = "27.898.47"
$ws.Range("E20").Value = "  +2.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.86"
$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.127"
$ws.Range("E22").Value = "  +0.87%  "

$ws.Range("D23").Value = "2.127.08"
$ws.Range("E23").Value = "  +1.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.86"
$ws.Range("E24").Value = "  +1.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.947"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.02"
$ws.Range("E26").Value = "  +2.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.58"
$ws.Range("E27").Value = "  +0.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.057"
$ws.Range("E28").Value = "  +3.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.03"
$ws.Range("E29").Value = "  +0.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.994"
$ws.Range("E30").Value = "  +1.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08924"
$ws.Range("E31").Value = "  +0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.331"
$ws.Range("E32").Value = "  +0.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.230"
$ws.Range("E33").Value = "  +4.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7696"
$ws.Range("E34").Value = "  +4.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.660"
$ws.Range("E35").Value = "  +2.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.734"
$ws.Range("E36").Value = "  +1.69%  "

$ws.Range("E37").Value = "  +1.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02050"
$ws.Range("E38").Value = "  +2.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5609"
$ws.Range("E39").Value = "  +5.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05369"
$ws.Range("E40").Value = "  +2.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.996"
$ws.Range("E41").Value = "  +0.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.051"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.583"
$ws.Range("E43").Value = "  +2.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1535"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4885"
$ws.Range("E45").Value = "  +2.70%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.68"
$ws.Range("E46").Value = "  +0.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.46"
$ws.Range("E47").Value = "  +3.00%  "

$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.670"
$ws.Range("E49").Value = "  +2.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.14"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06115"
$ws.Range("E51").Value = "  +0.79%  "

